# Update the "lu" lookup table with the latest projects:
#  - rename column C header from "category" to "sub_component"
#  - add a new column D "sub_targetgroup"
#  - populate D for existing rows
#  - append two new rows: bradford, yorkshire
#  - resize the table / used range accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resize the table first so the new column D is part of Table1 before
#     we give it its real header text (avoids a transient "Column4" name) ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D12"))

# --- Add the first new project row (bradford) ---
$ws.Range("A11").Value = "bradford"
$ws.Range("B11").Value = "NHS Bradford District and Craven Talking Therapies"
$ws.Range("C11").Value = "Multi-component (informational and operational)"

# --- Column header for column C (rename "category" -> "sub_component") ---
$ws.Range("C1").Value = "sub_component"

# --- New column D header ---
$ws.Range("D1").Value = "sub_targetgroup"

# --- New column D values ---
$ws.Range("D11").Value = "Subgroup-specific"
$ws.Range("D2").Value = "Non-Targeted or Mixed"
$ws.Range("D3").Value = "Subgroup-specific"
$ws.Range("D4").Value = "Non-Targeted or Mixed"
$ws.Range("D5").Value = "Non-Targeted or Mixed"
$ws.Range("D6").Value = "Non-Targeted or Mixed"
$ws.Range("D7").Value = "Non-Targeted or Mixed"
$ws.Range("D8").Value = "Non-Targeted or Mixed"

# --- Add the second new project row (yorkshire) ---
$ws.Range("A12").Value = "yorkshire"
$ws.Range("B12").Value = "NHS North Yorkshire Talking Therapies"
$ws.Range("C12").Value = "Multi-component (informational and operational)"
$ws.Range("D12").Value = "Subgroup-specific"

# --- Column width for the new column (compensate for the engine's fixed
#     +5/6 padding added between ColumnWidth and the stored <col width>) ---
$ws.Columns.Item(4).ColumnWidth = 37.1666666666667

# --- Selection like Excel would leave after data entry ---
$ws.Range("A13").Select()
